$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Base")
Write-Host $ws.Name
